$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.353.03'
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").Value = '2.246.44'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.67%  '

$ws.Range("E8").Value = '  +1.52%  '

$ws.Range("D9").Value = '2.290.15'
$ws.Range("E9").Value = '  +1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0949'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.78%  '

$ws.Range("E11").Value = '  +2.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.325'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.84%  '

$ws.Range("D14").Value = '2.671.67'
$ws.Range("E14").Value = '  +0.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.78'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.60%  '

$ws.Range("D16").Value = '54.176.64'
$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("E17").Value = '  +0.76%  '

$ws.Range("D18").Value = '2.272.18'
$ws.Range("E18").Value = '  +0.53%  '

$ws.Range("E19").Value = '  +4.71%  '

$ws.Range("E20").Value = '  +3.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '302.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.49%  '

$ws.Range("E24").Value = '  -1.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.23%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("E27").Value = '  +2.40%  '

$ws.Range("D28").Value = '2.404.54'
$ws.Range("E28").Value = '  +0.93%  '

$ws.Range("E29").Value = '  +5.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.43%  '

$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").Value = '0.0₃0690'
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("E34").Value = '  +2.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.63'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("E39").Value = '  +3.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.870'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.77%  '

$ws.Range("E42").Value = '  -0.97%  '

$ws.Range("E43").Value = '  +3.06%  '

$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("E45").Value = '  +1.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '128.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0891'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.545'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '240.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0486'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.87%  '
